$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-06-15 Thursday", $true, $true, $false, $false, $false, $true, 1, $false, "2023-06-16 Friday", 2) | Out-Null
$d.Content.Find.Execute("76-64=", $true, $true, $false, $false, $false, $true, 1, $false, "71-63=", 2) | Out-Null
$d.Content.Find.Execute("54+9=", $true, $true, $false, $false, $false, $true, 1, $false, "57+5=", 2) | Out-Null
$d.Content.Find.Execute("95-24=", $true, $true, $false, $false, $false, $true, 1, $false, "91-38=", 2) | Out-Null
$d.Content.Find.Execute("39+53=", $true, $true, $false, $false, $false, $true, 1, $false, "28+7=", 2) | Out-Null
$d.Content.Find.Execute("91-60=", $true, $true, $false, $false, $false, $true, 1, $false, "39-20=", 2) | Out-Null
$d.Content.Find.Execute("69-5=", $true, $true, $false, $false, $false, $true, 1, $false, "73+7=", 2) | Out-Null
$d.Content.Find.Execute("53-3=", $true, $true, $false, $false, $false, $true, 1, $false, "64-53=", 2) | Out-Null
$d.Content.Find.Execute("41+40=", $true, $true, $false, $false, $false, $true, 1, $false, "48-1=", 2) | Out-Null
$d.Content.Find.Execute("20+11=", $true, $true, $false, $false, $false, $true, 1, $false, "83+9=", 2) | Out-Null
$d.Content.Find.Execute("10+79=", $true, $true, $false, $false, $false, $true, 1, $false, "73-56=", 2) | Out-Null
$d.Content.Find.Execute("24+40=", $true, $true, $false, $false, $false, $true, 1, $false, "6+82=", 2) | Out-Null
$d.Content.Find.Execute("78-12=", $true, $true, $false, $false, $false, $true, 1, $false, "38-36=", 2) | Out-Null
$d.Content.Find.Execute("23+37=", $true, $true, $false, $false, $false, $true, 1, $false, "58-28=", 2) | Out-Null
$d.Content.Find.Execute("72-60=", $true, $true, $false, $false, $false, $true, 1, $false, "44+45=", 2) | Out-Null
$d.Content.Find.Execute("59-56=", $true, $true, $false, $false, $false, $true, 1, $false, "18+57=", 2) | Out-Null
$d.Content.Find.Execute("54-45=", $true, $true, $false, $false, $false, $true, 1, $false, "24+31=", 2) | Out-Null
$d.Content.Find.Execute("70+6=", $true, $true, $false, $false, $false, $true, 1, $false, "53+1=", 2) | Out-Null
$d.Content.Find.Execute("48-14=", $true, $true, $false, $false, $false, $true, 1, $false, "69-18=", 2) | Out-Null
$d.Content.Find.Execute("52+21=", $true, $true, $false, $false, $false, $true, 1, $false, "55+0=", 2) | Out-Null
$d.Content.Find.Execute("49+1=", $true, $true, $false, $false, $false, $true, 1, $false, "4+8=", 2) | Out-Null
$d.Content.Find.Execute("26-1=", $true, $true, $false, $false, $false, $true, 1, $false, "45+47=", 2) | Out-Null
$d.Content.Find.Execute("89-51=", $true, $true, $false, $false, $false, $true, 1, $false, "85-23=", 2) | Out-Null
$d.Content.Find.Execute("98-12=", $true, $true, $false, $false, $false, $true, 1, $false, "89-42=", 2) | Out-Null
$d.Content.Find.Execute("60+16=", $true, $true, $false, $false, $false, $true, 1, $false, "61-36=", 2) | Out-Null
$d.Content.Find.Execute("4+71=", $true, $true, $false, $false, $false, $true, 1, $false, "42-22=", 2) | Out-Null
$d.Content.Find.Execute("63+34=", $true, $true, $false, $false, $false, $true, 1, $false, "80-19=", 2) | Out-Null
$d.Content.Find.Execute("26+16=", $true, $true, $false, $false, $false, $true, 1, $false, "18+40=", 2) | Out-Null
$d.Content.Find.Execute("56-35=", $true, $true, $false, $false, $false, $true, 1, $false, "69+6=", 2) | Out-Null
$d.Content.Find.Execute("55-19=", $true, $true, $false, $false, $false, $true, 1, $false, "62+25=", 2) | Out-Null
$d.Content.Find.Execute("16+4=", $true, $true, $false, $false, $false, $true, 1, $false, "31+2=", 2) | Out-Null
$d.Content.Find.Execute("11+39=", $true, $true, $false, $false, $false, $true, 1, $false, "52-35=", 2) | Out-Null
$d.Content.Find.Execute("5+57=", $true, $true, $false, $false, $false, $true, 1, $false, "67-33=", 2) | Out-Null
$d.Content.Find.Execute("90-28=", $true, $true, $false, $false, $false, $true, 1, $false, "46-8=", 2) | Out-Null
$d.Content.Find.Execute("76-13=", $true, $true, $false, $false, $false, $true, 1, $false, "80-37=", 2) | Out-Null
$d.Content.Find.Execute("53-29=", $true, $true, $false, $false, $false, $true, 1, $false, "5+24=", 2) | Out-Null
$d.Content.Find.Execute("74-66=", $true, $true, $false, $false, $false, $true, 1, $false, "3+15=", 2) | Out-Null
$d.Content.Find.Execute("42+47=", $true, $true, $false, $false, $false, $true, 1, $false, "9+31=", 2) | Out-Null
$d.Content.Find.Execute("12+4=", $true, $true, $false, $false, $false, $true, 1, $false, "40+9=", 2) | Out-Null
$d.Content.Find.Execute("76+5=", $true, $true, $false, $false, $false, $true, 1, $false, "1+67=", 2) | Out-Null
$d.Content.Find.Execute("98-16=", $true, $true, $false, $false, $false, $true, 1, $false, "70+26=", 2) | Out-Null
$d.Content.Find.Execute("31+47=", $true, $true, $false, $false, $false, $true, 1, $false, "41-9=", 2) | Out-Null
$d.Content.Find.Execute("42-41=", $true, $true, $false, $false, $false, $true, 1, $false, "90+9=", 2) | Out-Null
$d.Content.Find.Execute("89-50=", $true, $true, $false, $false, $false, $true, 1, $false, "2+87=", 2) | Out-Null
$d.Content.Find.Execute("68-60=", $true, $true, $false, $false, $false, $true, 1, $false, "94-41=", 2) | Out-Null
$d.Content.Find.Execute("94-62=", $true, $true, $false, $false, $false, $true, 1, $false, "98-54=", 2) | Out-Null
$d.Content.Find.Execute("14+16=", $true, $true, $false, $false, $false, $true, 1, $false, "30+63=", 2) | Out-Null
$d.Content.Find.Execute("12+58=", $true, $true, $false, $false, $false, $true, 1, $false, "10+43=", 2) | Out-Null
$d.Content.Find.Execute("40+56=", $true, $true, $false, $false, $false, $true, 1, $false, "82-23=", 2) | Out-Null
$d.Content.Find.Execute("12+42=", $true, $true, $false, $false, $false, $true, 1, $false, "10-2=", 2) | Out-Null
$d.Content.Find.Execute("89-83=", $true, $true, $false, $false, $false, $true, 1, $false, "16+76=", 2) | Out-Null
$d.Content.Find.Execute("47+29=", $true, $true, $false, $false, $false, $true, 1, $false, "53-45=", 2) | Out-Null
$d.Content.Find.Execute("76-66=", $true, $true, $false, $false, $false, $true, 1, $false, "20+72=", 2) | Out-Null
$d.Content.Find.Execute("16-15=", $true, $true, $false, $false, $false, $true, 1, $false, "33+55=", 2) | Out-Null
$d.Content.Find.Execute("5+50=", $true, $true, $false, $false, $false, $true, 1, $false, "37+20=", 2) | Out-Null
$d.Content.Find.Execute("69-65=", $true, $true, $false, $false, $false, $true, 1, $false, "79-7=", 2) | Out-Null
$d.Content.Find.Execute("55-46=", $true, $true, $false, $false, $false, $true, 1, $false, "24+60=", 2) | Out-Null
$d.Content.Find.Execute("44+43=", $true, $true, $false, $false, $false, $true, 1, $false, "91-9=", 2) | Out-Null
$d.Content.Find.Execute("99-34=", $true, $true, $false, $false, $false, $true, 1, $false, "41+14=", 2) | Out-Null
$d.Content.Find.Execute("29+20=", $true, $true, $false, $false, $false, $true, 1, $false, "4+78=", 2) | Out-Null
$d.Content.Find.Execute("58+16=", $true, $true, $false, $false, $false, $true, 1, $false, "51+24=", 2) | Out-Null
$d.Content.Find.Execute("92-43=", $true, $true, $false, $false, $false, $true, 1, $false, "71-45=", 2) | Out-Null
$d.Content.Find.Execute("79-31=", $true, $true, $false, $false, $false, $true, 1, $false, "50+17=", 2) | Out-Null
$d.Content.Find.Execute("8+85=", $true, $true, $false, $false, $false, $true, 1, $false, "7+56=", 2) | Out-Null
$d.Content.Find.Execute("10+65=", $true, $true, $false, $false, $false, $true, 1, $false, "20+2=", 2) | Out-Null
$d.Content.Find.Execute("50-31=", $true, $true, $false, $false, $false, $true, 1, $false, "4+20=", 2) | Out-Null
$d.Content.Find.Execute("6+74=", $true, $true, $false, $false, $false, $true, 1, $false, "75-63=", 2) | Out-Null
$d.Content.Find.Execute("39+12=", $true, $true, $false, $false, $false, $true, 1, $false, "1+44=", 2) | Out-Null
$d.Content.Find.Execute("29+61=", $true, $true, $false, $false, $false, $true, 1, $false, "37-23=", 2) | Out-Null
$d.Content.Find.Execute("35-14=", $true, $true, $false, $false, $false, $true, 1, $false, "26+50=", 2) | Out-Null
$d.Content.Find.Execute("47-35=", $true, $true, $false, $false, $false, $true, 1, $false, "60-55=", 2) | Out-Null
$d.Content.Find.Execute("5+75=", $true, $true, $false, $false, $false, $true, 1, $false, "74+16=", 2) | Out-Null
$d.Content.Find.Execute("78-75=", $true, $true, $false, $false, $false, $true, 1, $false, "99-11=", 2) | Out-Null
$d.Content.Find.Execute("15+8=", $true, $true, $false, $false, $false, $true, 1, $false, "49-18=", 2) | Out-Null
$d.Content.Find.Execute("91-78=", $true, $true, $false, $false, $false, $true, 1, $false, "81+15=", 2) | Out-Null
$d.Content.Find.Execute("13+14=", $true, $true, $false, $false, $false, $true, 1, $false, "47+7=", 2) | Out-Null
$d.Content.Find.Execute("11+53=", $true, $true, $false, $false, $false, $true, 1, $false, "26-14=", 2) | Out-Null
$d.Content.Find.Execute("3+28=", $true, $true, $false, $false, $false, $true, 1, $false, "31-17=", 2) | Out-Null
$d.Content.Find.Execute("89-89=", $true, $true, $false, $false, $false, $true, 1, $false, "2+41=", 2) | Out-Null
$d.Content.Find.Execute("72-61=", $true, $true, $false, $false, $false, $true, 1, $false, "7+23=", 2) | Out-Null
$d.Content.Find.Execute("9+42=", $true, $true, $false, $false, $false, $true, 1, $false, "20+79=", 2) | Out-Null
$d.Content.Find.Execute("35+37=", $true, $true, $false, $false, $false, $true, 1, $false, "35+38=", 2) | Out-Null
$d.Content.Find.Execute("48+41=", $true, $true, $false, $false, $false, $true, 1, $false, "29+34=", 2) | Out-Null
$d.Content.Find.Execute("32+21=", $true, $true, $false, $false, $false, $true, 1, $false, "74-10=", 2) | Out-Null
$d.Content.Find.Execute("71-19=", $true, $true, $false, $false, $false, $true, 1, $false, "60-49=", 2) | Out-Null
$d.Content.Find.Execute("28+27=", $true, $true, $false, $false, $false, $true, 1, $false, "37-34=", 2) | Out-Null
$d.Content.Find.Execute("92-58=", $true, $true, $false, $false, $false, $true, 1, $false, "58-6=", 2) | Out-Null
$d.Content.Find.Execute("48-30=", $true, $true, $false, $false, $false, $true, 1, $false, "78-0=", 2) | Out-Null
$d.Content.Find.Execute("30+68=", $true, $true, $false, $false, $false, $true, 1, $false, "75-29=", 2) | Out-Null
$d.Content.Find.Execute("96-6=", $true, $true, $false, $false, $false, $true, 1, $false, "64-17=", 2) | Out-Null
$d.Content.Find.Execute("32+17=", $true, $true, $false, $false, $false, $true, 1, $false, "96-73=", 2) | Out-Null
$d.Content.Find.Execute("81+2=", $true, $true, $false, $false, $false, $true, 1, $false, "94-55=", 2) | Out-Null
$d.Content.Find.Execute("69-3=", $true, $true, $false, $false, $false, $true, 1, $false, "3+22=", 2) | Out-Null
$d.Content.Find.Execute("4+45=", $true, $true, $false, $false, $false, $true, 1, $false, "6+56=", 2) | Out-Null
$d.Content.Find.Execute("49+4=", $true, $true, $false, $false, $false, $true, 1, $false, "3+22=", 2) | Out-Null
$d.Content.Find.Execute("79+5=", $true, $true, $false, $false, $false, $true, 1, $false, "2+34=", 2) | Out-Null
$d.Content.Find.Execute("49+32=", $true, $true, $false, $false, $false, $true, 1, $false, "5+60=", 2) | Out-Null
$d.Content.Find.Execute("15+28=", $true, $true, $false, $false, $false, $true, 1, $false, "60-43=", 2) | Out-Null
$d.Content.Find.Execute("22+18=", $true, $true, $false, $false, $false, $true, 1, $false, "57-37=", 2) | Out-Null
$d.Content.Find.Execute("19-7=", $true, $true, $false, $false, $false, $true, 1, $false, "78-59=", 2) | Out-Null
$d.Content.Find.Execute("58-37=", $true, $true, $false, $false, $false, $true, 1, $false, "60+13=", 2) | Out-Null
